$d = $word.ActiveDocument

# Before this edit the document ends with:
#   ...paragraph 5: "2022年6月3日星期五"
#   paragraph 6 (last): "中雨，今天是农历五月初五，中国传统端午节了；" +
#                        "端午节，这一天我们要吃粽子，赛龙舟。"   (two runs)
#
# After this edit it should end with:
#   paragraph 5 (unchanged): "2022年6月3日星期五"
#   paragraph 6 (new):       "中雨，今天是农历五月初五，中国传统端午节了；"
#   paragraph 7 (new):       "2022" + "年6月7日星期二"
#   paragraph 8 (was 6):     "晴，今天是高考的第一天，上午考语文，下午考数学。"

$p5 = $d.Paragraphs.Item(5)
$lastPara = $d.Paragraphs.Last
$lastStart = $lastPara.Range.Start

# --- Step 1: build the new "中雨..." paragraph right after paragraph 5.
# Grab the formatted text of the last paragraph's first run (it already has the
# exact text + rFonts hint=eastAsia formatting we need) and insert it just before
# paragraph 5's own paragraph mark, then split it off into its own paragraph.
# Using Range.FormattedText (instead of typing fresh text) copies the run's exact
# formatting, and inserting+splitting this way (rather than
# InsertParagraphAfter()-then-fill) avoids leaving a stray empty run behind.
$run1Text = "中雨，今天是农历五月初五，中国传统端午节了；"
$srcRun1 = $d.Range($lastStart, $lastStart + $run1Text.Length)
$ft1 = $srcRun1.FormattedText

$insertAt1 = $p5.Range.End - 1
$insertPoint1 = $d.Range($insertAt1, $insertAt1)
$insertPoint1.FormattedText = $ft1
$d.Range($insertAt1, $insertAt1).InsertParagraphAfter()

# --- Step 2: build the new "2022年6月7日星期二" paragraph right after that one.
# Paragraph 5 already has exactly the run split we need ("2022" plain +
# "年6月3日星期五" hint=eastAsia), so duplicate its formatted text the same way
# (excluding its trailing paragraph mark!), then fix up the date text with a
# Find/Replace scoped to the new paragraph only.
$p6 = $d.Paragraphs.Item(6)
$srcP5Text = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$ft2 = $srcP5Text.FormattedText

$insertAt2 = $p6.Range.End - 1
$insertPoint2 = $d.Range($insertAt2, $insertAt2)
$insertPoint2.FormattedText = $ft2
$d.Range($insertAt2, $insertAt2).InsertParagraphAfter()

$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute("年6月3日星期五", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "年6月7日星期二", 2)

# --- Step 3: the original last paragraph (now paragraph 8) replaces its two runs
# with a single new sentence. Replace the text only (paragraph mark excluded), so
# the paragraph's own pPr/formatting is left untouched.
$last = $d.Paragraphs.Last
$body = $d.Range($last.Range.Start, $last.Range.End - 1)
$body.Text = "晴，今天是高考的第一天，上午考语文，下午考数学。"
